$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $text) {
    $scratch = $sheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $sheet.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
    $scratch.ClearFormats()
}

Set-CellText $ws "D2" "41.528.77"
Set-CellText $ws "E2" "  +0.73%  "
Set-CellText $ws "D3" "2.479.66"
Set-CellText $ws "E3" "  +0.56%  "
Set-CellText $ws "D4" "0.998"
Set-CellText $ws "E4" "  -0.14%  "
Set-CellText $ws "D5" "313.51"
Set-CellText $ws "E5" "  +0.41%  "
Set-CellText $ws "D6" "92.90"
Set-CellText $ws "E6" "  -1.28%  "
Set-CellText $ws "E7" "  -1.20%  "
Set-CellText $ws "D8" "0.999"
Set-CellText $ws "E8" "  -0.21%  "
Set-CellText $ws "D9" "0.507"
Set-CellText $ws "E9" "  +1.89%  "
Set-CellText $ws "D10" "32.72"
Set-CellText $ws "E10" "  -1.86%  "
Set-CellText $ws "D11" "0.0784"
Set-CellText $ws "E11" "  +0.79%  "
Set-CellText $ws "E12" "  +1.98%  "
Set-CellText $ws "D13" "2.860.92"
Set-CellText $ws "E13" "  +0.56%  "
Set-CellText $ws "B14" "Polkadot"
Set-CellText $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText $ws "D14" "6.85"
Set-CellText $ws "E14" "  -1.90%  "
Set-CellText $ws "B15" "Chainlink"
Set-CellText $ws "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-CellText $ws "D15" "16.18"
Set-CellText $ws "E15" "  +9.20%  "
Set-CellText $ws "D16" "2.464.00"
Set-CellText $ws "E16" "  -0.45%  "
Set-CellText $ws "D17" "0.767"
Set-CellText $ws "E17" "  -2.04%  "
Set-CellText $ws "D18" "41.515.00"
Set-CellText $ws "E18" "  +0.80%  "
Set-CellText $ws "D19" "6.40"
Set-CellText $ws "E19" "  +1.68%  "
Set-CellText $ws "D20" "0.0₃0939"
Set-CellText $ws "E20" "  +2.01%  "
Set-CellText $ws "D21" "71.78"
Set-CellText $ws "E21" "  +5.02%  "
Set-CellText $ws "D22" "11.29"
Set-CellText $ws "E22" "  -0.12%  "
Set-CellText $ws "D23" "236.76"
Set-CellText $ws "E23" "  +0.62%  "
Set-CellText $ws "D24" "2.70"
Set-CellText $ws "E25" "  -0.10%  "
Set-CellText $ws "E26" "  -0.50%  "
Set-CellText $ws "D27" "24.93"
Set-CellText $ws "E27" "  +4.20%  "
Set-CellText $ws "E28" "  +0.03%  "
Set-CellText $ws "D29" "9.64"
Set-CellText $ws "E29" "  +0.38%  "
Set-CellText $ws "D30" "35.94"
Set-CellText $ws "E30" "  -1.44%  "
Set-CellText $ws "D31" "158.22"
Set-CellText $ws "E31" "  +3.56%  "
Set-CellText $ws "D32" "5.45"
Set-CellText $ws "E32" "  -0.59%  "
Set-CellText $ws "E33" "  +0.96%  "
Set-CellText $ws "D34" "0.0757"
Set-CellText $ws "E34" "  +1.70%  "
Set-CellText $ws "E35" "  -8.34%  "
Set-CellText $ws "E36" "  +1.83%  "
Set-CellText $ws "D37" "0.106"
Set-CellText $ws "E37" "  +3.37%  "
Set-CellText $ws "D38" "2.91"
Set-CellText $ws "E38" "  -4.20%  "
Set-CellText $ws "D39" "1.84"
Set-CellText $ws "E39" "  -2.62%  "
Set-CellText $ws "E40" "  -0.06%  "
Set-CellText $ws "D41" "4.10"
Set-CellText $ws "E41" "  -2.63%  "
Set-CellText $ws "E42" "  -0.26%  "
Set-CellText $ws "D43" "19.39"
Set-CellText $ws "E43" "  -3.80%  "
Set-CellText $ws "D44" "1.979.24"
Set-CellText $ws "E44" "  +0.58%  "
Set-CellText $ws "E45" "  -0.13%  "
Set-CellText $ws "D46" "2.95"
Set-CellText $ws "E46" "  -2.76%  "
Set-CellText $ws "D47" "8.97"
Set-CellText $ws "E47" "  +3.44%  "
Set-CellText $ws "D48" "2.719.93"
Set-CellText $ws "E48" "  +0.52%  "
Set-CellText $ws "D49" "97.70"
Set-CellText $ws "E49" "  +0.53%  "
Set-CellText $ws "D50" "68.13"
Set-CellText $ws "E50" "  -1.71%  "
Set-CellText $ws "D51" "72.29"
Set-CellText $ws "E51" "  -2.08%  "
